$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '69.428.33'
Set-TextValue "E2" '  +2.45%  '

Set-TextValue "D3" '3.391.44'
Set-TextValue "E3" '  +1.65%  '

Set-TextValue "E4" '  +0.01%  '

Set-TextValue "D5" '587.34'
Set-TextValue "E5" '  +1.17%  '

Set-TextValue "D6" '180.25'
Set-TextValue "E6" '  +2.60%  '

Set-TextValue "E7" '  -0.04%  '

Set-TextValue "D8" '0.597'
Set-TextValue "E8" '  +1.32%  '

Set-TextValue "D9" '0.197'
Set-TextValue "E9" '  +8.09%  '

Set-TextValue "E10" '  +1.94%  '

Set-TextValue "D11" '48.59'
Set-TextValue "E11" '  +3.65%  '

Set-TextValue "D12" '0.0000283'
Set-TextValue "E12" '  +4.09%  '

Set-TextValue "D13" '679.50'
Set-TextValue "E13" '  -1.66%  '

Set-TextValue "D14" '8.66'
Set-TextValue "E14" '  +2.48%  '

Set-TextValue "D15" '3.931.51'
Set-TextValue "E15" '  +1.47%  '

Set-TextValue "D16" '69.454.85'
Set-TextValue "E16" '  +2.44%  '

Set-TextValue "D17" '0.121'
Set-TextValue "E17" '  +1.67%  '

Set-TextValue "D18" '3.373.50'
Set-TextValue "E18" '  +0.95%  '

Set-TextValue "D19" '17.73'
Set-TextValue "E19" '  +1.08%  '

Set-TextValue "D20" '11.30'
Set-TextValue "E20" '  +2.25%  '

Set-TextValue "D21" '0.906'
Set-TextValue "E21" '  +1.40%  '

Set-TextValue "D22" '5.44'
Set-TextValue "E22" '  +0.35%  '

Set-TextValue "D23" '17.15'
Set-TextValue "E23" '  +1.18%  '

Set-TextValue "D24" '103.26'
Set-TextValue "E24" '  +2.45%  '

Set-TextValue "D25" '3.93'
Set-TextValue "E25" '  +0.71%  '

Set-TextValue "D26" '2.74'
Set-TextValue "E26" '  +1.72%  '

Set-TextValue "D27" '9.65'
Set-TextValue "E27" '  +1.26%  '

Set-TextValue "D28" '33.98'
Set-TextValue "E28" '  +2.88%  '

Set-TextValue "D29" '8.76'
Set-TextValue "E29" '  +2.29%  '

Set-TextValue "D30" '6.97'
Set-TextValue "E30" '  -1.25%  '

Set-TextValue "D31" '11.15'
Set-TextValue "E31" '  +1.41%  '

Set-TextValue "D32" '558.58'
Set-TextValue "E32" '  -1.55%  '

Set-TextValue "E33" '  +1.09%  '

Set-TextValue "B34" 'dogwifhat'
Set-TextValue "C34" 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D34" '3.58'
Set-TextValue "E34" '  +8.12%  '

Set-TextValue "B35" 'OKB'
Set-TextValue "C35" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D35" '58.61'
Set-TextValue "E35" '  +2.18%  '

Set-TextValue "E36" '  +0.13%  '

Set-TextValue "D37" '3.678.47'
Set-TextValue "E37" '  -0.72%  '

Set-TextValue "D38" '35.91'
Set-TextValue "E38" '  +1.95%  '

Set-TextValue "E39" '  +4.17%  '

Set-TextValue "B40" 'PEPE'
Set-TextValue "C40" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue "D40" '0.0₃0724'
Set-TextValue "E40" '  +7.61%  '

Set-TextValue "B41" 'Stacks'
Set-TextValue "C41" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D41" '3.28'
Set-TextValue "E41" '  +3.58%  '

Set-TextValue "B42" 'Fetch.AI'
Set-TextValue "C42" 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue "D42" '2.69'
Set-TextValue "E42" '  +2.40%  '

Set-TextValue "D43" '0.341'
Set-TextValue "E43" '  +1.66%  '

Set-TextValue "B44" 'ApeXProtocol'
Set-TextValue "C44" 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue "D44" '3.32'
Set-TextValue "E44" '  +1.02%  '

Set-TextValue "B45" 'VeChain'
Set-TextValue "C45" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D45" '0.0423'
Set-TextValue "E45" '  +3.75%  '

Set-TextValue "E46" '  +0.74%  '

Set-TextValue "E47" '  +1.23%  '

Set-TextValue "E48" '  +6.46%  '

Set-TextValue "E49" '  -0.05%  '

Set-TextValue "D50" '133.44'
Set-TextValue "E50" '  +1.14%  '

Set-TextValue "E51" '  +3.45%  '
